# daily auto push: 2026-01-21 02:31 UTC
# Insert two new 30-minute-ranking rows for 2026/01/21 (time slots 8 and 9)
# right after the existing 2026/01/21 rows (680-681), pushing every
# subsequent row down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the old row 682 ("2026/12/29" block).
# Everything currently at row 682 onward (through 723) shifts down to
# 684-725, carrying its original values/format with it.
$ws.Rows("682:683").Insert()

# The new rows should hold text dates like the rest of column A, not an
# auto-converted date serial, so force text formatting before assigning.
$ws.Range("A682:A683").NumberFormat = "@"

$ws.Range("A682").Value = "2026/01/21"
$ws.Range("B682").Value = "水"
$ws.Range("C682").Value = 8
$ws.Range("D682").Value = 201

$ws.Range("A683").Value = "2026/01/21"
$ws.Range("B683").Value = "水"
$ws.Range("C683").Value = 9
$ws.Range("D683").Value = 201
